$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update group info: GroupName1 -> GroupName2, TestGroup1 -> TestGroup2
$ws.Range("A1").Value = "GroupName2"
$ws.Range("A2").Value = "TestGroup2"
$ws.Range("A3").Value = "notes2"

# Move the selection to A3
$ws.Range("A3").Select()
